$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize throughput values (divide by 1000) for rows 2-7, columns B:G ---
$data = @{
    2 = @(20.3699645691919, 38.8322600344195, 66.8705963839594, 106.833823024382, 111.919291653341, 126.83998550265)
    3 = @(11.0902928008241, 21.2847129835871, 36.0428821545353, 63.4505200661696, 77.2609967437496, 85.0554998633008)
    4 = @(7.39717293699421, 14.362207472883899, 27.044928950241, 43.2852792035124, 66.0693739823969, 71.2966462094782)
    5 = @(5.56760920844406, 10.7762160248169, 20.6453358896836, 34.553053405, 53.0762245981195, 61.5303036962158)
    6 = @(4.50407345439706, 8.69516559712687, 16.5045035486076, 28.5900564671765, 46.5260875344787, 54.2624907942168)
    7 = @(3.517884538, 6.669522001, 12.73520769, 23.19605455, 36.84479128, 46.4815669)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # column B = 2 ... G = 7
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# --- Remove the three blank spacer rows (10-12) entirely ---
$ws.Range("F10:G12").Clear()

# --- Clear the B:G contents (and formatting) in rows 13-15, leaving only column A ---
$ws.Range("B13:G15").Clear()

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("C10").Select()
